$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 156.8
$ws.Range("I9").Value = 301.5
$ws.Range("K9").Value = 301.5
$ws.Range("M9").Value = -132.5
$ws.Range("H113").Value = 5615.75
$ws.Range("I113").Value = 2981.5
$ws.Range("K113").Value = 2981.5
$ws.Range("M113").Value = 272.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3625.4736
$ws.Range("I61").Value = 2535.2727
$ws.Range("K61").Value = 2535.2727
$ws.Range("M61").Value = -2323.2727
$ws.Range("H74").Value = 2283
$ws.Range("I74").Value = 2172.15
$ws.Range("J74").Value = 4500
$ws.Range("K74").Value = 2172.15
$ws.Range("L74").Value = 4500
$ws.Range("M74").Value = -1298.15
$ws.Range("N74").Value = -6248
$ws.Range("H77").Value = 2283
$ws.Range("I77").Value = 2172.15
$ws.Range("J77").Value = 4500
$ws.Range("K77").Value = 10860.75
$ws.Range("L77").Value = 22500
$ws.Range("M77").Value = -6492.75
$ws.Range("N77").Value = -31236
$ws.Range("H97").Value = 706.75
$ws.Range("I97").Value = 706.75
$ws.Range("K97").Value = 706.75
$ws.Range("M97").Value = -210.75
$ws.Range("H132").Value = 3662.5
$ws.Range("I132").Value = 3665.375
$ws.Range("J132").Value = 3639.5
$ws.Range("K132").Value = 10996.125
$ws.Range("L132").Value = 10918.5
$ws.Range("M132").Value = -8466.125
$ws.Range("N132").Value = -15978.5
$ws.Range("H136").Value = 3625.4736
$ws.Range("I136").Value = 2535.2727
$ws.Range("K136").Value = 7605.8181
$ws.Range("M136").Value = -5055.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1448.3
$ws.Range("I20").Value = 1309.625
$ws.Range("K20").Value = 1309.625
$ws.Range("M20").Value = -1062.625
$ws.Range("H94").Value = 1423.1364
$ws.Range("I94").Value = 1460.8948
$ws.Range("K94").Value = 1460.8948
$ws.Range("M94").Value = -1009.8948
$ws.Range("H107").Value = 5263.75
$ws.Range("I107").Value = 4945.6
$ws.Range("J107").Value = 5491
$ws.Range("K107").Value = 4945.6
$ws.Range("L107").Value = 5491
$ws.Range("M107").Value = -3025.6
$ws.Range("N107").Value = -9331
$ws.Range("H134").Value = 1570.7778
$ws.Range("I134").Value = 1570.7778
$ws.Range("K134").Value = 4712.3334
$ws.Range("M134").Value = -2177.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5587.853
$ws.Range("I31").Value = 796.7778
$ws.Range("J31").Value = 7312.64
$ws.Range("K31").Value = 796.7778
$ws.Range("L31").Value = 7312.64
$ws.Range("M31").Value = -501.7778
$ws.Range("N31").Value = -7902.64
$ws.Range("H34").Value = 5587.853
$ws.Range("I34").Value = 796.7778
$ws.Range("J34").Value = 7312.64
$ws.Range("K34").Value = 796.7778
$ws.Range("L34").Value = 7312.64
$ws.Range("M34").Value = -594.7778
$ws.Range("N34").Value = -7716.64
$ws.Range("H58").Value = 4259.7144
$ws.Range("I58").Value = 3706.25
$ws.Range("K58").Value = 3706.25
$ws.Range("M58").Value = -3503.25
$ws.Range("H132").Value = 3014.7058
$ws.Range("I132").Value = 2589.2856
$ws.Range("K132").Value = 7767.8568
$ws.Range("M132").Value = -5237.8568
$ws.Range("H134").Value = 1218
$ws.Range("I134").Value = 1372.5
$ws.Range("K134").Value = 4117.5
$ws.Range("M134").Value = -1582.5
$ws.Range("H136").Value = 4259.7144
$ws.Range("I136").Value = 3706.25
$ws.Range("K136").Value = 11118.75
$ws.Range("M136").Value = -8568.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 1919.75
$ws.Range("I134").Value = 1919.75
$ws.Range("K134").Value = 5759.25
$ws.Range("M134").Value = -689.25
$ws.Range("H138").Value = 4197.8
$ws.Range("I138").Value = 2998.6
$ws.Range("J138").Value = 5397
$ws.Range("K138").Value = 8995.799999999999
$ws.Range("L138").Value = 16191
$ws.Range("M138").Value = -3855.799999999999
$ws.Range("N138").Value = -26471
$ws.Range("H139").Value = 4340.5
$ws.Range("I139").Value = 3710.75
$ws.Range("J139").Value = 5600
$ws.Range("K139").Value = 11132.25
$ws.Range("L139").Value = 16800
$ws.Range("M139").Value = -5992.25
$ws.Range("N139").Value = -27080
$ws.Range("H140").Value = 3905.8
$ws.Range("I140").Value = 3632.25
$ws.Range("J140").Value = 5000
$ws.Range("K140").Value = 10896.75
$ws.Range("L140").Value = 15000
$ws.Range("M140").Value = -5716.75
$ws.Range("N140").Value = -25360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 2500
$ws.Range("I21").Value = 2500
$ws.Range("K21").Value = 2500
$ws.Range("M21").Value = -2327
$ws.Range("H30").Value = 2500
$ws.Range("I30").Value = 2500
$ws.Range("K30").Value = 2500
$ws.Range("M30").Value = -2395
$ws.Range("H122").Value = 1596.5
$ws.Range("I122").Value = 1596.5
$ws.Range("K122").Value = 4789.5
$ws.Range("M122").Value = -2339.5
$ws.Range("H126").Value = 4250
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 4666.6665
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 13999.9995
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -18939.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 333666
$ws.Range("I20").Value = 499.5
$ws.Range("J20").Value = 999999
$ws.Range("K20").Value = 499.5
$ws.Range("L20").Value = 999999
$ws.Range("M20").Value = -273.5
$ws.Range("N20").Value = -1000451
$ws.Range("H46").Value = 2914.3333
$ws.Range("I46").Value = 2371.5
$ws.Range("K46").Value = 2371.5
$ws.Range("M46").Value = -2183.5
$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 2000
$ws.Range("K132").Value = 6000
$ws.Range("M132").Value = -3470
$ws.Range("H136").Value = 3146
$ws.Range("I136").Value = 2075.75
$ws.Range("K136").Value = 6227.25
$ws.Range("M136").Value = -3677.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 23587
$ws.Range("J2").Value = 23587
$ws.Range("L2").Value = 23587
$ws.Range("N2").Value = -23811
$ws.Range("H4").Value = 27206.666
$ws.Range("I4").Value = 22648
$ws.Range("J4").Value = 50000
$ws.Range("K4").Value = 22648
$ws.Range("L4").Value = 50000
$ws.Range("M4").Value = -22535
$ws.Range("N4").Value = -50226
$ws.Range("H132").Value = 2468.5
$ws.Range("I132").Value = 2374.6667
$ws.Range("J132").Value = 2750
$ws.Range("K132").Value = 7124.000100000001
$ws.Range("L132").Value = 8250
$ws.Range("M132").Value = -4594.000100000001
$ws.Range("N132").Value = -13310
$ws.Range("H136").Value = 3838.7083
$ws.Range("I136").Value = 2992.7334
$ws.Range("J136").Value = 5248.6665
$ws.Range("K136").Value = 8978.200199999999
$ws.Range("L136").Value = 15745.9995
$ws.Range("M136").Value = -6428.200199999999
$ws.Range("N136").Value = -20845.9995
